# Automatische test-sync: 2025-08-19 21:04:50
# Append the new mail-log entry (row 23) to the "Logs" sheet, extend the
# conditional-formatting ranges that cover the data rows, and bump the
# "Intern verzoek / Actie voor medewerker" tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")

# --- New row 23 data ---------------------------------------------------
$logs.Range("A23").Value = "Vraag over product"
$logs.Range("B23").Value = "documentatie@testbedrijf123.nl"
$logs.Range("D23").Value = "Intern verzoek / Actie voor medewerker"
$logs.Range("F23").Value = "2025-08-19 21:04:22"
$logs.Range("G23").Value = "Nee"
$logs.Range("H23").Value = "Ja"
$logs.Range("I23").Value = "Nee"
$logs.Range("J23").Value = "Nee"

# --- Extend conditional formatting ranges from row 22 to row 23 --------
$logs.Range("D2:D22").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D23"))
$logs.Range("G2:G22").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G23"))
$logs.Range("H2:H22").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H23"))
$logs.Range("I2:I22").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I23"))
$logs.Range("J2:J22").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J23"))

# --- Update Dashboard tally ---------------------------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 22
